$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

# ---------------------------------------------------------------------------
# The sheet currently ends with:
#   row 59 -> data row (item 53)
#   row 60 -> totals row (only P60/Q60 populated)
#   row 61 -> footer row (date / page / developed-by)
#
# We need to insert a new data row (item 54) before the totals row, which
# pushes the totals row to 61 (with an updated total) and the footer row to
# 62. We do this with targeted Copy/PasteSpecial(Formats) + explicit value
# assignment (cell by cell) instead of Rows.Insert(), because a whole-row
# Insert duplicates style records in this engine; per-cell format copy
# reuses the existing style indices.
# ---------------------------------------------------------------------------

# Step 1: move the footer row (61) down to row 62, formats + values.
foreach ($col in $cols) {
    $src = $ws.Range($col + "61")
    $dst = $ws.Range($col + "62")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    if ($src.Value2 -ne "") {
        $dst.Value = $src.Value2
    }
}
$ws.Rows(62).RowHeight = 16.5

# Step 2: move the totals row (60) down to row 61, formats + values.
$ws.Range("P60").Copy()
$ws.Range("P61").PasteSpecial(-4122)
$ws.Range("P61").Value = $ws.Range("P60").Value2

$ws.Range("Q60").Copy()
$ws.Range("Q61").PasteSpecial(-4122)

$ws.Rows(61).RowHeight = 25.5

# The rest of row 61 (A:O) is unused by the totals row - clear any leftover
# footer formatting/values that used to live at row 61 before the shift.
$ws.Range("A61:O61").Value = ""
$ws.Range("A61:O61").Style = "Normal"

# Update the grand total to include the new item's sell price (30.00).
$ws.Range("P61").Value = 3216.125

# Step 3: build the new data row 60 (item 54) by copying formats from the
# row above (59, item 53) and filling in the new item's values.
foreach ($col in $cols) {
    $src = $ws.Range($col + "59")
    $dst = $ws.Range($col + "60")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}
$ws.Rows(60).RowHeight = 24.75

$ws.Range("A60").Value = 54
$ws.Range("C60").Value = "معجون اسنان فلورو بالكولا"
$ws.Range("H60").Value = "4:0"
$ws.Range("L60").Value = "0"
$ws.Range("N60").Value = "30.00"
$ws.Range("P60").Value = "30.0000"
$ws.Range("Q60").Value = "1:0"

# ---------------------------------------------------------------------------
# Merged cells: recreate the per-row merges for the shifted rows 61/62 and
# add the new merges for row 60.
# ---------------------------------------------------------------------------
$ws.Range("P61:Q61").Merge()
$ws.Range("A62:F62").Merge()
$ws.Range("G62:I62").Merge()
$ws.Range("K62:Q62").Merge()

$ws.Range("A60:B60").Merge()
$ws.Range("C60:G60").Merge()
$ws.Range("H60:K60").Merge()
$ws.Range("L60:M60").Merge()
$ws.Range("N60:O60").Merge()
